$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "Here is a summary of the code:`n" +
    "**Component:** ``HeroCarousel```n" +
    "**Purpose:** Display a carousel of images with hero product images`n" +
    "**Key Features:**`n" +
    "1. Uses React and Next.js`n" +
    "2. Imports Carousel component from react-responsive-carousel library`n" +
    "3. Defines an array of image objects (``heroImages``) with src, alt, and other attributes`n" +
    "4. Renders the carousel with auto-play and infinite loop features`n" +
    "5. Displays a custom arrow icon above the carousel using Image component`n" +
    "6. Maps through ``heroImages`` array to render each image in the carousel`n" +
    "**Notable CSS classes:**`n" +
    "1. ``.max-xl:hidden``: Hides the element on large screens`n" +
    "2. ``object-contain``: Maintains aspect ratio of images within the carousel"

$ws.Range("B10").Value = $newText
